$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Way to make the gate for small female <al>lizards</al>"
#   Move the space that currently begins the (Arial-styled) " small female "
#   run so that it instead ends the (plain-styled) " the gate for" run.
#   Net effect on rendered text: none ("...for small..." unchanged); only the
#   run boundary / formatting of that one space character changes.
# ---------------------------------------------------------------------------
$r1 = $d.Content
$null = $r1.Find.Execute("the gate for small")
# "the gate for" is 12 characters; the boundary we want sits right after it,
# i.e. immediately before the existing space character.
$boundary1 = $r1.Start + 12
$insPos1 = $d.Range($boundary1, $boundary1)
$insPos1.InsertAfter(" ")
# The original leading space of " small female " has now been pushed one
# character later (it immediately follows the space we just inserted) -
# delete it so the visible text is unchanged.
$oldSpace1 = $d.Range($boundary1 + 1, $boundary1 + 2)
$oldSpace1.Delete()

# ---------------------------------------------------------------------------
# Edit 2: "other," -> "others,"
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "other, that are applied in the same way", $true, $false, $false, $false,
    $false, $true, 1, $false, "others, that are applied in the same way", 2)

# ---------------------------------------------------------------------------
# Edit 3: "But take heed to make sure that the end of the " ->
#         "But take heed to make sure, that with the end of the "
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "But take heed to make sure that the end of the ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "But take heed to make sure, that with the end of the ", 2)

# ---------------------------------------------------------------------------
# Edit 4a: insert a comma right after the styled "</tl>" run (keeping that
#          run's Courier-New / blue formatting) so it reads "</tl>,".
# ---------------------------------------------------------------------------
$r4 = $d.Content
$null = $r4.Find.Execute("hot iron wire</tl>")
$tlEnd = $r4.End
$insPos4 = $d.Range($tlEnd, $tlEnd)
$insPos4.InsertAfter(",")

# Edit 4b: rewrite the plain run that used to read
#   "  hardly touches the animal, for the sand of the second gate will not
#    touch this part.  But make sure that the end of the "
# into
#   " the end of the wax barely touches the animal, for the sand of the
#    second cast will not touch this part. But make sure that the end of the "
$null = $d.Content.Find.Execute(
    "  hardly touches the animal, for the sand of the second gate will not touch this part.  But make sure that the end of the ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " the end of the wax barely touches the animal, for the sand of the second cast will not touch this part. But make sure that the end of the ",
    2)

# ---------------------------------------------------------------------------
# Edit 5: " they will serve as <fr>abreuvoirs</fr> for the molded thing."
#   The literal "<fr>"/"</fr>" markup (currently plain text inside one run)
#   needs to become its own Courier-New / blue run, same as the many other
#   <tag> markers elsewhere in the document.
# ---------------------------------------------------------------------------
# 5a. Strip the literal tag text back down to bare "abreuvoirs".
$null = $d.Content.Find.Execute(
    "as <fr>abreuvoirs</fr> for", $true, $false, $false, $false, $false,
    $true, 1, $false, "as abreuvoirs for", 2)

# 5b. Borrow the formatting of an existing "<fr>"/"</fr>" pair elsewhere in
#     the document (e.g. around "forge") via FormattedText, and splice
#     (read-only) copies of it in immediately before/after "abreuvoirs".
$donor5 = $d.Content
$null = $donor5.Find.Execute("<fr>forge</fr>")
$donorBase5 = $donor5.Start
$openDonor5 = $d.Range($donorBase5, $donorBase5 + 4)
$closeDonor5 = $d.Range($donorBase5 + 9, $donorBase5 + 14)
$openFT5 = $openDonor5.FormattedText
$closeFT5 = $closeDonor5.FormattedText

# Insert the closing tag first (later position) so the earlier insertion
# point used afterwards is unaffected by this one's length change.
$tgtClose5 = $d.Content
$null = $tgtClose5.Find.Execute("abreuvoirs")
$afterPos5 = $d.Range($tgtClose5.End, $tgtClose5.End)
$afterPos5.FormattedText = $closeFT5

$tgtOpen5 = $d.Content
$null = $tgtOpen5.Find.Execute("abreuvoirs")
$beforePos5 = $d.Range($tgtOpen5.Start, $tgtOpen5.Start)
$beforePos5.FormattedText = $openFT5

# ---------------------------------------------------------------------------
# Edit 6: "The most important is that the ears of the " ->
#         "The principal thing is that the ears of the "
#   "principal thing" must land in a new, plainly-formatted run (matching
#   the style of ordinary body text elsewhere), while "The " and
#   " is that the ears of the " retain their explicit Arial formatting.
# ---------------------------------------------------------------------------
# Borrow plain formatting from the "hot iron wire" run: temporarily retext
# it to get a FormattedText carrying plain formatting + the new words, use
# that to replace "most important", then restore the donor run's own text.
$donor6 = $d.Content
$null = $donor6.Find.Execute("hot iron wire")
$donorFT6 = $donor6.FormattedText
$donorFT6.Text = "principal thing"

$tgt6 = $d.Content
$null = $tgt6.Find.Execute("most important")
$tgt6.FormattedText = $donorFT6

$restore6 = $d.Content
$null = $restore6.Find.Execute("principal thing")
$restore6.Text = "hot iron wire"
